# Auto-update draw results: append the 2025-10-13 Pick 4 draw as a new row
# at the bottom of the results table (mirrors the nightly scraper job that
# appends one row per day, all columns stored as plain text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 27

# The existing rows store every value as literal text (dates like
# "2025-10-13", zero-padded phase codes like "251013", etc. are NOT real
# numbers/dates). Pre-format the new row as Text so Excel's "smart" entry
# parsing doesn't convert the date-looking / number-looking strings into a
# real date serial or numeric value.
$newRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$newRange.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-10-13"
$ws.Range("B" + $newRow).Value = "Pick 4"
$ws.Range("C" + $newRow).Value = "251013"
$ws.Range("D" + $newRow).Value = "3-8-8-5"
$ws.Range("E" + $newRow).Value = "2025-10-13T21:37:37.024+04:00"

# Drop the explicit Text number-format again so the new row doesn't end up
# styled differently from the rows above it (they rely on the sheet's
# default/general cell style even though their content is text).
$newRange.ClearFormats()

# These are numeric-looking strings on purpose; tell Excel's error checker
# not to flag them, same as the rest of the column already does.
for ($c = 1; $c -le 5; $c++) {
    try {
        $ws.Cells.Item($newRow, $c).Errors.Item(9).Ignore = $true
    } catch {
    }
}
